$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ItICM")

# Insert a new column before L (old L..AK shift right to M..AL).
$ws.Columns("L").Insert()

# Old column K held "ISIC 20T21: Chemicals and pharmaceutical products" -- it is
# split into two separate ISIC categories: the (still-in-place) K column becomes
# "ISIC 20: Chemicals" and the newly inserted L column becomes
# "ISIC 21: Pharmaceuticals".
$ws.Range("L1").Value = "ISIC 21: Pharmaceuticals"
$ws.Range("K1").Value = "ISIC 20: Chemicals"

# Populate the new column's data cells (defaults to 0/unmapped for every
# industry except "other industries", which also maps to pharmaceuticals).
$ws.Range("L2").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("L9").Value = 1
